$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The weekly data rows (2-8) get their Fecha/Volumen/Precio values re-shuffled
# across rows (columns D, J, K, L, M, P), while all other columns stay put.

$values = @{
    2 = @{ D = 44382; J = 160; K = 7000; L = 8000; M = 7438; P = 124 }
    3 = @{ D = 44281; J = 120; K = 5500; L = 6000; M = 5750; P = 96 }
    4 = @{ D = 44362; J = 120; K = 8000; L = 9000; M = 8500; P = 142 }
    5 = @{ D = 44242; J = 160; K = 5000; L = 5500; M = 5250; P = 88 }
    6 = @{ D = 44421; J = 100; K = 8000; L = 9000; M = 8500; P = 142 }
    7 = @{ D = 44400; J = 120; K = 9000; L = 10000; M = 9500; P = 158 }
    8 = @{ D = 44494; J = 120; K = 5000; L = 6000; M = 5500; P = 92 }
}

foreach ($row in $values.Keys) {
    $rowVals = $values[$row]
    $ws.Range("D$row").Value = $rowVals.D
    $ws.Range("J$row").Value = $rowVals.J
    $ws.Range("K$row").Value = $rowVals.K
    $ws.Range("L$row").Value = $rowVals.L
    $ws.Range("M$row").Value = $rowVals.M
    $ws.Range("P$row").Value = $rowVals.P
}
